# Rename the "wide" (WInd energy moDEls) taxonomy prefix/namespace/title
# to "idem" (IDEM: wInD Energy Models), matching the regenerated .ttl
# export from the Google Sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = 'http://purl.org/idem/'
$ws.Range("B3").Value = 'idem'
$ws.Range("C3").Value = 'http://purl.org/idem/'
$ws.Range("A17").Value = 'idem:MeteorologicalModels'
$ws.Range("A18").Value = 'idem:OceanographicModels'
$ws.Range("A19").Value = 'idem:HydrodynamicModels'
$ws.Range("F19").Value = 'idem:OceanographicModels'
$ws.Range("A20").Value = 'idem:WaveModels'
$ws.Range("F20").Value = 'idem:OceanographicModels'
$ws.Range("A21").Value = 'idem:FlowModels'
$ws.Range("A22").Value = 'idem:AerolasticModels'
$ws.Range("A23").Value = 'idem:ElectricalModels'
$ws.Range("A24").Value = 'idem:FinancialModels'
$ws.Range("A25").Value = 'idem:GCM'
$ws.Range("F25").Value = 'idem:MeteorologicalModels'
$ws.Range("A26").Value = 'idem:Mesoscale'
$ws.Range("F26").Value = 'idem:MeteorologicalModels'
$ws.Range("A27").Value = 'idem:Hindcast'
$ws.Range("F27").Value = 'idem:MeteorologicalModels'
$ws.Range("A28").Value = 'idem:MorisonEquation'
$ws.Range("F28").Value = 'idem:HydrodynamicModels'
$ws.Range("A29").Value = 'idem:Radiation-Diffraction'
$ws.Range("F29").Value = 'idem:HydrodynamicModels'
$ws.Range("A30").Value = 'idem:Linear'
$ws.Range("F30").Value = 'idem:WaveModels'
$ws.Range("A31").Value = 'idem:SecondOrder'
$ws.Range("F31").Value = 'idem:WaveModels'
$ws.Range("A32").Value = 'idem:FullyNonlinear'
$ws.Range("F32").Value = 'idem:WaveModels'
$ws.Range("A33").Value = 'idem:Linearized'
$ws.Range("F33").Value = 'idem:FlowModels'
$ws.Range("A34").Value = 'idem:RANS'
$ws.Range("F34").Value = 'idem:FlowModels'
$ws.Range("A35").Value = 'idem:LES'
$ws.Range("F35").Value = 'idem:FlowModels'
$ws.Range("A36").Value = 'idem:DNS'
$ws.Range("F36").Value = 'idem:FlowModels'
$ws.Range("A37").Value = 'idem:VortexMethod'
$ws.Range("F37").Value = 'idem:FlowModels'
$ws.Range("A38").Value = 'idem:Analytical'
$ws.Range("F38").Value = 'idem:FlowModels'
$ws.Range("A39").Value = 'idem:Physical'
$ws.Range("F39").Value = 'idem:FlowModels'
$ws.Range("A40").Value = 'idem:FEM'
$ws.Range("F40").Value = 'idem:AerolasticModels'
$ws.Range("A41").Value = 'idem:PowerFlow'
$ws.Range("F41").Value = 'idem:ElectricalModels'
$ws.Range("A42").Value = 'idem:OPF'
$ws.Range("F42").Value = 'idem:ElectricalModels'
$ws.Range("A43").Value = 'idem:Small-SignalModels'
$ws.Range("F43").Value = 'idem:ElectricalModels'
$ws.Range("A44").Value = 'idem:DynamicModels'
$ws.Range("F44").Value = 'idem:ElectricalModels'
$ws.Range("A45").Value = 'idem:ShortCircuitModels'
$ws.Range("F45").Value = 'idem:ElectricalModels'
$ws.Range("A46").Value = 'idem:StateEstimation'
$ws.Range("F46").Value = 'idem:ElectricalModels'
$ws.Range("A47").Value = 'idem:PowerProtectionAnalysisModels'
$ws.Range("F47").Value = 'idem:ElectricalModels'
$ws.Range("A48").Value = 'idem:ContingencyAnalysisModels'
$ws.Range("F48").Value = 'idem:ElectricalModels'
$ws.Range("A49").Value = 'idem:HarmonicModels'
$ws.Range("F49").Value = 'idem:ElectricalModels'
$ws.Range("A50").Value = 'idem:BoS'
$ws.Range("F50").Value = 'idem:FinancialModels'
$ws.Range("A51").Value = 'idem:NPV'
$ws.Range("F51").Value = 'idem:FinancialModels'
$ws.Range("A52").Value = 'idem:LCOE'
$ws.Range("F52").Value = 'idem:FinancialModels'
$ws.Range("A53").Value = 'idem:IRR'
$ws.Range("F53").Value = 'idem:FinancialModels'
$ws.Range("A54").Value = 'idem:'

$ws.Range("B8").Value = "IDEM: wInD Energy Models`n"
